$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "end" marker from C4 up to C3, giving it the same
# "Input" named style used by the other key cells (e.g. A2/B3).
$ws.Range("C3").Value = "end"
$ws.Range("C3").Style = "Input"

# B4 and C4 revert to the default "Check Cell" style, and C4's
# old "end" text is cleared out since it now lives in C3.
$ws.Range("B4").Style = "Check Cell"
$ws.Range("C4").ClearContents()
$ws.Range("C4").Style = "Check Cell"

# Update the active selection to E4 to match the saved view state.
$ws.Range("E4").Select()
